# Re-ran projection with updated TNMs and using 3 allocation periods
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# CMAP Region updated net migration projections
$ws.Range("C3").Value = 115000
$ws.Range("C4").Value = 200000
$ws.Range("C5").Value = 225000
$ws.Range("C6").Value = 250000
$ws.Range("C7").Value = 225000

# External IL updated net migration projections
$ws.Range("C8").Value = 15000
$ws.Range("C11").Value = 30000
$ws.Range("C12").Value = 25000

# External IN updated net migration projections
$ws.Range("C19").Value = 15000

# External WI updated net migration projections
$ws.Range("C22").Value = 20000
$ws.Range("C23").Value = 15000
$ws.Range("C24").Value = 10000

# Update the active selection to reflect where the author last worked
$ws.Range("H13").Select()
